# Duplicate row 3 ("Assured Clear Distance Ahead" / 4511.21(A)) into a new
# row 4, pushing all subsequent rows down by one (rows 4-40 become 5-41).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("3:3").Copy()
$ws.Rows("4:4").Insert()

# Match the author's final selection on the newly inserted row.
$ws.Range("A4:D4").Select()
